$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A student is being disconnected from the group -> append their record as a
# new row right after the last one (row 7), mirroring the layout/format of
# the existing roster rows (e.g. row 5).
$ws.Range("A8").Value = 146824
$ws.Range("B8").Value = "arturo"
$ws.Range("C8").Value = "ledezma"
$ws.Range("D8").Value = "macias"
$ws.Range("E8").Value = "M"
$ws.Range("F8").Value = "25/01/2010"
$ws.Range("G8").Value = "er@ugto.mx"
$ws.Range("H8").Value = "LISC"

# Match the number formatting used by the existing roster rows so the new
# entry renders the same way (row 5 is a same-shaped, fully populated row;
# the "M"/sex column follows the header-row style instead).
$ws.Range("A8:D8").NumberFormat = $ws.Range("A5:D5").NumberFormat
$ws.Range("E8").NumberFormat = $ws.Range("E3").NumberFormat
$ws.Range("F8").NumberFormat = $ws.Range("F5").NumberFormat
$ws.Range("G8:H8").NumberFormat = $ws.Range("G5:H5").NumberFormat

$ws.Range("H8").Select()
